# Auto-generated update script for violent crime workbook (adds 2022-05-04 data)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 2049
$ws.Range("H3").Value = 8349
$ws.Range("I3").Value = 2161
$ws.Range("H4").Value = 1658
$ws.Range("I4").Value = 542
$ws.Range("I5").Value = 192
$ws.Range("I6").Value = 2582
$ws.Range("H7").Value = 25969
$ws.Range("I7").Value = 7526

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 77

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I3").Value = 9
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I3").Value = 71
$ws.Range("I6").Value = 68
$ws.Range("I7").Value = 240

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 49
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 102
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 289

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 51
$ws.Range("I5").Value = 8
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 78
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 255
$ws.Range("I8").Value = 476
$ws.Range("I14").Value = 40
$ws.Range("I15").Value = 93
$ws.Range("I19").Value = 212
$ws.Range("I20").Value = 202
$ws.Range("I27").Value = 71
$ws.Range("I29").Value = 486
$ws.Range("I36").Value = 96
$ws.Range("I37").Value = 240
$ws.Range("I42").Value = 252
$ws.Range("I48").Value = 76
$ws.Range("I49").Value = 44
$ws.Range("I52").Value = 151
$ws.Range("I54").Value = 176
$ws.Range("I60").Value = 42
$ws.Range("I61").Value = 8
$ws.Range("H63").Value = 190
$ws.Range("I63").Value = 38
$ws.Range("I65").Value = 179
$ws.Range("I67").Value = 289
$ws.Range("I71").Value = 18
$ws.Range("I72").Value = 27
$ws.Range("I78").Value = 99
$ws.Range("I79").Value = 195
$ws.Range("I80").Value = 24
$ws.Range("I85").Value = 354
$ws.Range("I87").Value = 8
$ws.Range("I89").Value = 77
$ws.Range("I90").Value = 87
$ws.Range("I97").Value = 58
$ws.Range("I98").Value = 52
$ws.Range("I99").Value = 140
$ws.Range("H101").Value = 25969
$ws.Range("I101").Value = 7526

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 88
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 150
$ws.Range("I3").Value = 161
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 486

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 212

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I2").Value = 10
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 86
$ws.Range("I5").Value = 12
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 354

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 47

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 65
$ws.Range("I6").Value = 65
$ws.Range("I7").Value = 252

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 55
$ws.Range("I3").Value = 56
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 202

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 27
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 36
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 146
$ws.Range("I3").Value = 129
$ws.Range("I7").Value = 476

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I4").Value = 10
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 42

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("I2").Value = 4
$ws.Range("I7").Value = 18

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I3").Value = 6
$ws.Range("I7").Value = 27

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 24

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 74
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 8

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 8
